$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2986.4688
$ws.Range("I62").Value = 2722.6
$ws.Range("J62").Value = 6944.5
$ws.Range("K62").Value = 2722.6
$ws.Range("L62").Value = 6944.5
$ws.Range("M62").Value = -2098.6
$ws.Range("N62").Value = -8192.5

$ws.Range("H65").Value = 2986.4688
$ws.Range("I65").Value = 2722.6
$ws.Range("J65").Value = 6944.5
$ws.Range("K65").Value = 13613
$ws.Range("L65").Value = 34722.5
$ws.Range("M65").Value = -10493
$ws.Range("N65").Value = -40962.5

$ws.Range("H108").Value = 0
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()

$ws.Range("H138").Value = 3929.1282
$ws.Range("I138").Value = 1777.3846
$ws.Range("J138").Value = 5005
$ws.Range("K138").Value = 5332.1538
$ws.Range("L138").Value = 15015
$ws.Range("M138").Value = -192.1538
$ws.Range("N138").Value = -25295

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5944.8315
$ws.Range("I32").Value = 5773.427
$ws.Range("K32").Value = 5773.427
$ws.Range("M32").Value = -5486.427

$ws.Range("H61").Value = 8289.666999999999
$ws.Range("I61").Value = 8508.305
$ws.Range("K61").Value = 8508.305
$ws.Range("M61").Value = -8296.305

$ws.Range("H74").Value = 4862.8965
$ws.Range("I74").Value = 6466.1177
$ws.Range("J74").Value = 2591.6667
$ws.Range("K74").Value = 6466.1177
$ws.Range("L74").Value = 2591.6667
$ws.Range("M74").Value = -5592.1177
$ws.Range("N74").Value = -4339.6667

$ws.Range("H77").Value = 4862.8965
$ws.Range("I77").Value = 6466.1177
$ws.Range("J77").Value = 2591.6667
$ws.Range("K77").Value = 32330.5885
$ws.Range("L77").Value = 12958.3335
$ws.Range("M77").Value = -27962.5885
$ws.Range("N77").Value = -21694.3335

$ws.Range("H102").Value = 14724.934
$ws.Range("I102").Value = 29895.908
$ws.Range("K102").Value = 29895.908
$ws.Range("M102").Value = -28273.908

$ws.Range("H132").Value = 6272.636
$ws.Range("I132").Value = 3000
$ws.Range("J132").Value = 6599.9
$ws.Range("K132").Value = 9000
$ws.Range("L132").Value = 19799.7
$ws.Range("M132").Value = -6470
$ws.Range("N132").Value = -24859.7

$ws.Range("H136").Value = 8289.666999999999
$ws.Range("I136").Value = 8508.305
$ws.Range("K136").Value = 25524.915
$ws.Range("M136").Value = -22974.915

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H47").Value = 205888.33
$ws.Range("J47").Value = 205888.33
$ws.Range("L47").Value = 205888.33
$ws.Range("N47").Value = -206928.33

$ws.Range("H134").Value = 17939.1
$ws.Range("I134").Value = 43133.332
$ws.Range("J134").Value = 7141.5713
$ws.Range("K134").Value = 129399.996
$ws.Range("L134").Value = 21424.7139
$ws.Range("M134").Value = -126864.996
$ws.Range("N134").Value = -26494.7139

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 251190
$ws.Range("I16").Value = 1655.5
$ws.Range("J16").Value = 500724.5
$ws.Range("K16").Value = 1655.5
$ws.Range("L16").Value = 500724.5
$ws.Range("M16").Value = -1368.5
$ws.Range("N16").Value = -501298.5

$ws.Range("H31").Value = 8955.849
$ws.Range("J31").Value = 5993.1665
$ws.Range("L31").Value = 5993.1665
$ws.Range("N31").Value = -6583.1665

$ws.Range("H34").Value = 8955.849
$ws.Range("J34").Value = 5993.1665
$ws.Range("L34").Value = 5993.1665
$ws.Range("N34").Value = -6397.1665

$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()

$ws.Range("H113").Value = 251190
$ws.Range("I113").Value = 1655.5
$ws.Range("J113").Value = 500724.5
$ws.Range("K113").Value = 1655.5
$ws.Range("L113").Value = 500724.5
$ws.Range("M113").Value = 514.5
$ws.Range("N113").Value = -505064.5

$ws.Range("H134").Value = 5502.788
$ws.Range("I134").Value = 5817.96
$ws.Range("K134").Value = 17453.88
$ws.Range("M134").Value = -14918.88

$ws.Range("H141").Value = 199992.42
$ws.Range("J141").Value = 216947.95
$ws.Range("L141").Value = 216947.95
$ws.Range("N141").Value = -227307.95

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1739.75
$ws.Range("I34").Value = 1566.8572
$ws.Range("J34").Value = 2950
$ws.Range("K34").Value = 4700.571599999999
$ws.Range("L34").Value = 8850
$ws.Range("M34").Value = -4616.571599999999
$ws.Range("N34").Value = -9018

$ws.Range("H107").Value = 1248.15
$ws.Range("J107").Value = 1248.15
$ws.Range("L107").Value = 3744.45
$ws.Range("N107").Value = -7584.450000000001

$ws.Range("H113").Value = 758.05884
$ws.Range("I113").Value = 570
$ws.Range("J113").Value = 860.63635
$ws.Range("K113").Value = 1710
$ws.Range("L113").Value = 2581.90905
$ws.Range("M113").Value = 460
$ws.Range("N113").Value = -6921.90905

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 12456.737
$ws.Range("I70").Value = 12422.556
$ws.Range("K70").Value = 12422.556
$ws.Range("M70").Value = -12152.556

$ws.Range("H73").Value = 12456.737
$ws.Range("I73").Value = 12422.556
$ws.Range("K73").Value = 12422.556
$ws.Range("M73").Value = -11486.556

$ws.Range("H122").Value = 12178.046
$ws.Range("I122").Value = 8564.125
$ws.Range("K122").Value = 25692.375
$ws.Range("M122").Value = -23242.375

$ws.Range("H123").Value = 18554.334
$ws.Range("J123").Value = 18554.334
$ws.Range("L123").Value = 18554.334
$ws.Range("N123").Value = -23454.334

$ws.Range("H126").Value = 16185.6
$ws.Range("I126").Value = 61483.5
$ws.Range("J126").Value = 9216.691999999999
$ws.Range("K126").Value = 184450.5
$ws.Range("L126").Value = 27650.076
$ws.Range("M126").Value = -181980.5
$ws.Range("N126").Value = -32590.076

$ws.Range("H132").Value = 4907.1904
$ws.Range("I132").Value = 4907.1904
$ws.Range("K132").Value = 14721.5712
$ws.Range("M132").Value = -12191.5712

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 882.6667
$ws.Range("J9").Value = 1000
$ws.Range("L9").Value = 1000
$ws.Range("N9").Value = -1448

$ws.Range("H22").Value = 729.0833
$ws.Range("I22").Value = 694.9
$ws.Range("J22").Value = 900
$ws.Range("K22").Value = 694.9
$ws.Range("L22").Value = 900
$ws.Range("M22").Value = -399.9
$ws.Range("N22").Value = -1490

$ws.Range("H27").Value = 729.0833
$ws.Range("I27").Value = 694.9
$ws.Range("J27").Value = 900
$ws.Range("K27").Value = 694.9
$ws.Range("L27").Value = 900
$ws.Range("M27").Value = -587.9
$ws.Range("N27").Value = -1114

$ws.Range("H32").Value = 450
$ws.Range("I32").Value = 450
$ws.Range("K32").Value = 450
$ws.Range("M32").Value = -133

$ws.Range("H46").Value = 2366.4783
$ws.Range("I46").Value = 516.3333
$ws.Range("J46").Value = 2644
$ws.Range("K46").Value = 516.3333
$ws.Range("L46").Value = 2644
$ws.Range("M46").Value = -328.3333
$ws.Range("N46").Value = -3020

$ws.Range("H100").Value = 5471.607
$ws.Range("I100").Value = 4916.421
$ws.Range("K100").Value = 4916.421
$ws.Range("M100").Value = -4375.421

$ws.Range("H132").Value = 934279.5
$ws.Range("J132").Value = 4382
$ws.Range("L132").Value = 13146
$ws.Range("N132").Value = -18206

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 8229.229499999999
$ws.Range("I132").Value = 9214.444
$ws.Range("K132").Value = 27643.332
$ws.Range("M132").Value = -25113.332
